$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ============================================================
# "Create Deal Code Update" - add CreateDeal_/ShareDeal_ test rows
# Cell values are written in the same order the shared-string
# table was originally built in, so that the regenerated
# xl/sharedStrings.xml gets identical string indices.
# ============================================================

# --- Row 34: Test ID + Description ---
$s34A = @"
CreateDeal_TC001
"@
$ws.Range("A34").Value = $s34A
$s34B = @"
Validate whehter Shipper user is able to add new Deal  on following conditions.
a) Launch application and login applcation as Shipper admin
b) Goto Deals and click on add icon
c) Set mandatory fileds and clik on next.
d.) Click on Deals and Drafts
e.) Check whether the new deal widget added in Drafts
"@
$ws.Range("B34").Value = $s34B

# --- Run Mode column (C2:C33): Yes -> NO ---
for ($r = 2; $r -le 33; $r++) {
    $ws.Cells.Item($r, 3).Value = "NO"
}

# --- Row 35: Test ID + Description ---
$s35A = @"
CreateDeal_TC002
"@
$ws.Range("A35").Value = $s35A
$s35B = @"
Validate whehter Shipper admin is able to add new Deal  on following conditions.
a) Launch application and login applcation as Shipper user
b) Goto Deals and click on add icon
c) Set mandatory fileds and clik on next.
d.) Click on Deals and Drafts
e.) Check whether the new deal widget added in Drafts
"@
$ws.Range("B35").Value = $s35B

# --- Row 36: Test ID + Description ---
$s36A = @"
CreateDeal_TC003
"@
$ws.Range("A36").Value = $s36A
$s36B = @"
Validate whehter carrier is able to add new Deal  on following conditions.
a) Launch application and login applcation as Shipper user
b) Goto Deals and click on add icon
c) Set mandatory fileds and clik on next.
d.) Click on Deals and Drafts
e.) Check whether the new deal widget added in Drafts
"@
$ws.Range("B36").Value = $s36B

# --- Row 37: Test ID + Description ---
$s37A = @"
ShareDeal_TC001
"@
$ws.Range("A37").Value = $s37A
$s37B = @"
Validate whehter Shipper user is able to share new Deal  on following conditions.
a) Login with valid user id and Password 
b) Click on Deals menu
c) Click on Add New Deal button from Opportunity tab
d) Select Dealname,EquipmentTypes & # of Loads and click Next button
e) Select one or more contacts(Carrier) in contact list screen and click share button
f) Login as Carrier user and check opportunity tab whether Deal is displayed.
"@
$ws.Range("B37").Value = $s37B

# --- Row 37: Expected Result ("Deal has been shared") ---
$s37D = @"
Deal has been shared
"@
$ws.Range("D37").Value = $s37D

# --- Row 38: Test ID + Description ---
$s38A = @"
ShareDeal_TC002
"@
$ws.Range("A38").Value = $s38A
$s38B = @"
Validate whehter Shipper admin is able to share all new Deal  on following conditions.
a) Login withvalid user id and Password.
b) Click on Deals menu
c) Click on Add New Deal button from Opportunity tab
d) Select Dealname,EquipmentTypes & # of Loads and click Next button
e) Goto Drafts and check Deal widget available.
f) Select Deal, click on three dots and click on share.
g) Set ON Share All and click Share.
h) Login as Carrier user and check opportunity tab whether Deal is displayed.
"@
$ws.Range("B38").Value = $s38B

# --- Row 34: Expected Result ("Deal saved to draft") ---
$s34D = @"
Deal saved to draft
"@
$ws.Range("D34").Value = $s34D

# --- Remaining Expected Result cells (reuse strings minted above) ---
$ws.Range("D35").Value = $s34D
$ws.Range("D36").Value = $s34D
$ws.Range("D38").Value = $s37D

# --- Run Mode for the new rows stays "Yes" ---
for ($r = 34; $r -le 38; $r++) {
    $ws.Cells.Item($r, 3).Value = "Yes"
}

# --- Row heights for the new rows ---
$ws.Rows.Item(34).RowHeight = 105
$ws.Rows.Item(35).RowHeight = 105
$ws.Rows.Item(36).RowHeight = 105
$ws.Rows.Item(37).RowHeight = 165
$ws.Rows.Item(38).RowHeight = 180

# --- Formatting: column B wraps text; columns A, C, D are vertically centered ---
for ($r = 34; $r -le 38; $r++) {
    $ws.Cells.Item($r, 1).VerticalAlignment = -4108
    $ws.Cells.Item($r, 2).WrapText = $true
    $ws.Cells.Item($r, 3).VerticalAlignment = -4108
    $ws.Cells.Item($r, 4).VerticalAlignment = -4108
}

# --- Extend the AutoFilter + _FilterDatabase defined name over the new rows ---
$ws.Range("A1:F33").AutoFilter() | Out-Null
$ws.Range("A1:F38").AutoFilter() | Out-Null
$fdb = $wb.Names.Item("_xlnm._FilterDatabase")
$fdb.RefersTo = "='Automation Tests'!`$A`$1:`$F`$38"

# --- Move the selection/view down to the newly added rows ---
$ws.Range("B35").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 33
